# Update the customer id test data values on the "testdata" worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("testdata")

$ws.Range("B14").Value = "cus_N9wDkXhr3jDaG4"
$ws.Range("B18").Value = "cus_N9wDQ2uSnGUTFx"
